# Orégano - Mercado Mayorista Lo Valledor de Santiago
# Add a new weekly observation: insert a new data row at row 221 (just
# after the header + existing rows 2-220), pushing the remaining rows
# (old 221-300) down by one to become rows 222-301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 221; this shifts rows 221..300 to 222..301
# and extends the used range to A1:R301 automatically.
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with the new weekly record.
$ws.Range("A221").Value = 6
$ws.Range("B221").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C221").Value = "Metropolitana"
$ws.Range("D221").Value = 44988
$ws.Range("E221").Value = 13
$ws.Range("F221").Value = 100112029
$ws.Range("G221").Value = "Orégano"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 44
$ws.Range("K221").Value = 16000
$ws.Range("L221").Value = 17000
$ws.Range("M221").Value = 16455
$ws.Range("N221").Value = "$/docena de atados"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 5485
$ws.Range("Q221").Value = 3
$ws.Range("R221").Value = "Hortaliza"
